$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: update handoff/handback datetimes (rows 2 and 3 share same values)
$wsZh.Range("E2").Value = "2016-03-23 10:17:26"
$wsZh.Range("E3").Value = "2016-03-23 10:17:26"
$wsZh.Range("H2").Value = "2016-03-23 10:18:27"
$wsZh.Range("H3").Value = "2016-03-23 10:18:27"

# de-de sheet: update handoff/handback datetimes (rows 2 and 3 share same values)
$wsDe.Range("E2").Value = "2016-03-23 10:17:33"
$wsDe.Range("E3").Value = "2016-03-23 10:17:33"
$wsDe.Range("H2").Value = "2016-03-23 10:18:41"
$wsDe.Range("H3").Value = "2016-03-23 10:18:41"
